$d = $word.ActiveDocument

# 1. Normal style: add spacing after paragraph = 6pt (120 twips)
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 6

# 2. Add two new custom paragraph styles based on Heading4, following Normal,
#    mirroring the existing tei_* style family (teicollation, teiextent).
$collation = $d.Styles.Add("teicollation", 1)
$collation.NameLocal = "tei_collation"
$collation.BaseStyle = $d.Styles("Heading4")
$collation.NextParagraphStyle = $d.Styles("Normal")
$collation.QuickStyle = $true

$extent = $d.Styles.Add("teiextent", 1)
$extent.NameLocal = "tei_extent"
$extent.BaseStyle = $d.Styles("Heading4")
$extent.NextParagraphStyle = $d.Styles("Normal")
$extent.QuickStyle = $true
